$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D7: "Yes" -> "Resub"
$ws.Range("D7").Value = "Resub"

# E7: copy the date style/format from E6 (which already uses the date number format)
# then set the due-date value (2016-04-28, serial 42488)
$ws.Range("E6").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = 42488

# F7: new note text
$ws.Range("F7").Value = "Just D2"

# Update the active selection to F7, matching the saved view state
$ws.Range("F7").Select()
